$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "69.704.02"
Set-TextValue "E2" "  +0.25%  "

Set-TextValue "D3" "3.504.62"
Set-TextValue "E3" "  +0.26%  "

Set-TextValue "E4" "  -0.13%  "

Set-TextValue "D5" "603.96"
Set-TextValue "E5" "  -1.40%  "

Set-TextValue "D6" "195.53"
Set-TextValue "E6" "  +4.52%  "

Set-TextValue "E7" "  +0.23%  "

Set-TextValue "E8" "  +0.01%  "

Set-TextValue "D9" "0.201"
Set-TextValue "E9" "  -5.81%  "

Set-TextValue "D10" "0.650"
Set-TextValue "E10" "  +0.49%  "

Set-TextValue "D11" "53.47"
Set-TextValue "E11" "  +0.81%  "

Set-TextValue "D12" "0.0000302"
Set-TextValue "E12" "  -2.06%  "

Set-TextValue "D13" "9.49"
Set-TextValue "E13" "  -0.21%  "

Set-TextValue "D14" "4.062.36"
Set-TextValue "E14" "  +0.32%  "

Set-TextValue "D15" "594.49"
Set-TextValue "E15" "  -1.25%  "

Set-TextValue "D16" "12.81"
Set-TextValue "E16" "  +1.54%  "

Set-TextValue "D17" "69.822.92"
Set-TextValue "E17" "  +0.36%  "

Set-TextValue "D18" "19.04"
Set-TextValue "E18" "  +0.74%  "

Set-TextValue "E19" "  +1.78%  "

Set-TextValue "D20" "3.497.85"
Set-TextValue "E20" "  +0.00%  "

Set-TextValue "D21" "0.989"
Set-TextValue "E21" "  +0.21%  "

Set-TextValue "E23" "  +3.69%  "

Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "102.18"
Set-TextValue "E24" "  -3.33%  "

Set-TextValue "B25" "PancakeSwap"
Set-TextValue "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "4.65"
Set-TextValue "E25" "  +0.06%  "

Set-TextValue "D26" "3.16"
Set-TextValue "E26" "  +3.42%  "

Set-TextValue "D27" "10.84"
Set-TextValue "E27" "  -0.88%  "

Set-TextValue "D28" "9.56"
Set-TextValue "E28" "  -1.79%  "

Set-TextValue "D29" "33.30"
Set-TextValue "E29" "  -0.53%  "

Set-TextValue "D30" "4.33"
Set-TextValue "E30" "  +3.30%  "

Set-TextValue "D31" "7.04"
Set-TextValue "E31" "  +1.25%  "

Set-TextValue "D32" "12.41"
Set-TextValue "E32" "  -0.27%  "

Set-TextValue "E33" "  +0.00%  "

Set-TextValue "D34" "63.11"
Set-TextValue "E34" "  -0.54%  "

Set-TextValue "D35" "0.0₃0831"
Set-TextValue "E35" "  +6.84%  "

Set-TextValue "D36" "3.717.47"
Set-TextValue "E36" "  +3.14%  "

Set-TextValue "D37" "3.09"
Set-TextValue "E37" "  -2.71%  "

Set-TextValue "D38" "1.00"
Set-TextValue "E38" "  +0.01%  "

Set-TextValue "D39" "3.64"
Set-TextValue "E39" "  -1.81%  "

Set-TextValue "E40" "  -1.07%  "

Set-TextValue "E41" "  -1.04%  "

Set-TextValue "B42" "Kaspa"
Set-TextValue "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.133"
Set-TextValue "E42" "  -3.14%  "

Set-TextValue "B43" "Bittensor"
Set-TextValue "C43" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D43" "472.61"
Set-TextValue "E43" "  -8.69%  "

Set-TextValue "D44" "0.0454"
Set-TextValue "E44" "  -1.81%  "

Set-TextValue "D45" "0.140"
Set-TextValue "E45" "  -1.55%  "

Set-TextValue "D46" "2.82"
Set-TextValue "E46" "  -4.02%  "

Set-TextValue "E47" "  -1.50%  "

Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  +0.28%  "

Set-TextValue "D49" "8.42"
Set-TextValue "E49" "  -4.17%  "

Set-TextValue "D50" "0.000245"
Set-TextValue "E50" "  +1.89%  "

Set-TextValue "E51" "  +10.27%  "
